$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2033898305084746
$ws.Range("C2").Value = 0.53954802259887
$ws.Range("J2").Value = 0.002824858757062147
$ws.Range("P2").Value = 0.1525423728813559
$ws.Range("S2").Value = 0.1016949152542373
$ws.Range("B3").Value = 0.005154639175257732
$ws.Range("C3").Value = 0.005154639175257732
$ws.Range("J3").Value = 0.02577319587628866
$ws.Range("P3").Value = 0.7422680412371134
$ws.Range("S3").Value = 0.2216494845360825
$ws.Range("J4").Value = 0.09090909090909091
$ws.Range("P4").Value = 0.7272727272727273
$ws.Range("S4").Value = 0.1818181818181818
$ws.Range("B6").Value = 0.05882352941176471
$ws.Range("D6").Value = 0.009803921568627451
$ws.Range("E6").Value = 0.004901960784313725
$ws.Range("F6").Value = 0.04901960784313725
$ws.Range("J6").Value = 0.3137254901960784
$ws.Range("O6").Value = 0.02941176470588235
$ws.Range("Q6").Value = 0.1617647058823529
$ws.Range("R6").Value = 0.05392156862745098
$ws.Range("S6").Value = 0.3186274509803921
$ws.Range("B7").Value = 0.1469194312796208
$ws.Range("D7").Value = 0.004739336492890996
$ws.Range("F7").Value = 0.03317535545023697
$ws.Range("J7").Value = 0.1279620853080569
$ws.Range("O7").Value = 0.01421800947867299
$ws.Range("Q7").Value = 0.1943127962085308
$ws.Range("R7").Value = 0.07109004739336493
$ws.Range("S7").Value = 0.4075829383886256
$ws.Range("B8").Value = 0.1071428571428571
$ws.Range("D8").Value = 0.01275510204081633
$ws.Range("E8").Value = 0.00510204081632653
$ws.Range("F8").Value = 0.04591836734693878
$ws.Range("J8").Value = 0.1147959183673469
$ws.Range("O8").Value = 0.02040816326530612
$ws.Range("Q8").Value = 0.2295918367346939
$ws.Range("R8").Value = 0.08928571428571429
$ws.Range("S8").Value = 0.375
$ws.Range("B9").Value = 0.1347150259067358
$ws.Range("D9").Value = 0.02072538860103627
$ws.Range("F9").Value = 0.06735751295336788
$ws.Range("J9").Value = 0.1295336787564767
$ws.Range("O9").Value = 0.01036269430051814
$ws.Range("Q9").Value = 0.1865284974093264
$ws.Range("R9").Value = 0.08290155440414508
$ws.Range("S9").Value = 0.3678756476683938
$ws.Range("B10").Value = 0.1315175097276265
$ws.Range("D10").Value = 0.02490272373540856
$ws.Range("E10").Value = 0.0007782101167315176
$ws.Range("F10").Value = 0.05603112840466926
$ws.Range("J10").Value = 0.1175097276264591
$ws.Range("O10").Value = 0.01245136186770428
$ws.Range("Q10").Value = 0.2093385214007782
$ws.Range("R10").Value = 0.07859922178988327
$ws.Range("S10").Value = 0.3688715953307393
$ws.Range("G11").Value = 0.1451612903225807
$ws.Range("J11").Value = 0.07419354838709677
$ws.Range("K11").Value = 0.2096774193548387
$ws.Range("L11").Value = 0.5451612903225806
$ws.Range("S11").Value = 0.02580645161290323
$ws.Range("G12").Value = 0.7348066298342542
$ws.Range("J12").Value = 0.1602209944751381
$ws.Range("K12").Value = 0.01657458563535912
$ws.Range("L12").Value = 0.06629834254143646
$ws.Range("S12").Value = 0.02209944751381215
$ws.Range("G13").Value = 0.8085106382978723
$ws.Range("J13").Value = 0.1914893617021277
$ws.Range("F15").Value = 0.02926829268292683
$ws.Range("H15").Value = 0.1268292682926829
$ws.Range("I15").Value = 0.06829268292682927
$ws.Range("J15").Value = 0.3463414634146341
$ws.Range("K15").Value = 0.06341463414634146
$ws.Range("M15").Value = 0.00975609756097561
$ws.Range("N15").Value = 0.004878048780487805
$ws.Range("O15").Value = 0.08780487804878048
$ws.Range("S15").Value = 0.2634146341463415
$ws.Range("F16").Value = 0.05286343612334802
$ws.Range("H16").Value = 0.1277533039647577
$ws.Range("I16").Value = 0.1277533039647577
$ws.Range("J16").Value = 0.4361233480176211
$ws.Range("K16").Value = 0.0881057268722467
$ws.Range("M16").Value = 0.03083700440528634
$ws.Range("O16").Value = 0.02643171806167401
$ws.Range("S16").Value = 0.1101321585903084
$ws.Range("F17").Value = 0.01943844492440605
$ws.Range("H17").Value = 0.142548596112311
$ws.Range("I17").Value = 0.0734341252699784
$ws.Range("J17").Value = 0.4514038876889849
$ws.Range("K17").Value = 0.09719222462203024
$ws.Range("M17").Value = 0.02591792656587473
$ws.Range("O17").Value = 0.05399568034557235
$ws.Range("S17").Value = 0.1360691144708423
$ws.Range("F18").Value = 0.01704545454545454
$ws.Range("H18").Value = 0.1590909090909091
$ws.Range("I18").Value = 0.06818181818181818
$ws.Range("J18").Value = 0.4602272727272727
$ws.Range("K18").Value = 0.1022727272727273
$ws.Range("M18").Value = 0.01136363636363636
$ws.Range("O18").Value = 0.07386363636363637
$ws.Range("S18").Value = 0.1079545454545455
$ws.Range("F19").Value = 0.02625102543068089
$ws.Range("H19").Value = 0.2001640689089418
$ws.Range("I19").Value = 0.08531583264971287
$ws.Range("J19").Value = 0.3748974569319114
$ws.Range("K19").Value = 0.1148482362592289
$ws.Range("M19").Value = 0.01968826907301066
$ws.Range("N19").Value = 0.0008203445447087777
$ws.Range("O19").Value = 0.06808859721082855
$ws.Range("S19").Value = 0.1099261689909762
